$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list values (Price and Volume(1h) columns)
# Force text format so numeric-looking strings (e.g. "1.00") are preserved as text
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "60.296.76"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -2.08%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.369.45"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -2.28%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "566.15"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -2.43%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "139.93"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -6.73%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.368.34"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.68%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -4.03%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -2.82%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -1.42%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.945.48"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -2.32%  "
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +1.05%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "27.94"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.13%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.371.80"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -2.30%  "
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -3.92%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "60.433.60"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -2.07%  "
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -2.31%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.79"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -4.13%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "8.97"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -5.64%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "384.72"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -1.29%  "
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -2.33%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "72.96"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.09%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.03%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -8.47%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.521.81"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -1.82%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -2.34%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.03%  "
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -5.43%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -4.43%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.12"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -2.57%  "
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -9.54%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "23.47"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -2.33%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.400.48"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -2.04%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.88"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -2.63%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "167.88"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.57%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.91"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -6.26%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -5.04%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0768"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -2.72%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "27.00"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.09%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.00%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.777"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -2.08%  "
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -2.24%  "
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -2.08%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.514.15"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -3.29%  "
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -4.87%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "23.11"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.31%  "
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -2.88%  "
